$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 79; this shifts the existing rows 79-174 down to 80-175
# and keeps their contents intact (matches the diff: old row 79 -> new row 80, etc.,
# and old row 174 -> new row 175).
$ws.Rows("79:79").Insert()

# Populate the newly inserted row 79 with its data (same "context" values as the
# surrounding rows for this product/market, a new date, and its own measurements).
$ws.Range("A79").Value = 6
$ws.Range("B79").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C79").Value = "Metropolitana"
$ws.Range("D79").Value = 44601
$ws.Range("E79").Value = 13
$ws.Range("F79").Value = 100112001
$ws.Range("G79").Value = "Berenjena"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 400
$ws.Range("K79").Value = 8000
$ws.Range("L79").Value = 9000
$ws.Range("M79").Value = 8425
$ws.Range("N79").Value = "$/caja 60 unidades"
$ws.Range("O79").Value = "Región Metropolitana"
$ws.Range("P79").Value = 140
$ws.Range("Q79").Value = 60
$ws.Range("R79").Value = "Hortaliza"
